$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date (column C) from 45702 to 45703 for all data rows (2-41)
for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45702) {
        $cell.Value2 = 45703
    }
}

# Shift the "Beteckning" (column A) and "Area (ha)" (column G) values for rows 39-41
# New order: row39 <- old row41, row40 <- old row39, row41 <- old row40
$oldA39 = $ws.Cells.Item(39, 1).Value2
$oldA40 = $ws.Cells.Item(40, 1).Value2
$oldA41 = $ws.Cells.Item(41, 1).Value2

$oldG39 = $ws.Cells.Item(39, 7).Value2
$oldG40 = $ws.Cells.Item(40, 7).Value2
$oldG41 = $ws.Cells.Item(41, 7).Value2

$ws.Cells.Item(39, 1).Value2 = $oldA41
$ws.Cells.Item(40, 1).Value2 = $oldA39
$ws.Cells.Item(41, 1).Value2 = $oldA40

$ws.Cells.Item(39, 7).Value2 = $oldG41
$ws.Cells.Item(40, 7).Value2 = $oldG39
$ws.Cells.Item(41, 7).Value2 = $oldG40
